$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.473.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.629.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.86%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.40%  "
$ws.Range("E9").Value = "  -4.11%  "
$ws.Range("E10").Value = "  -4.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.364"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.098.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.82%  "
$ws.Range("E14").Value = "  -4.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.366.94"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.06%  "
$ws.Range("E16").Value = "  -3.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.631.17"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.44%  "
$ws.Range("E19").Value = "  -4.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "340.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.36%  "
$ws.Range("E21").Value = "  -7.67%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").Value = "  -3.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.86%  "
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0836"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("E30").Value = "  -4.35%  "
$ws.Range("E31").Value = "  -4.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "161.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.50%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.56%  "
$ws.Range("E35").Value = "  -4.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "19.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.10%  "
$ws.Range("E37").Value = "  -3.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "337.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.02%  "
$ws.Range("E39").Value = "  -2.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.909"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.09%  "
$ws.Range("E41").Value = "  -3.91%  "
$ws.Range("E42").Value = "  -1.84%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.612"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("E48").Value = "  -6.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0963"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.77%  "
$ws.Range("E51").Value = "  -5.14%  "

Write-Output "Update complete"
